$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells per the diff: rows 5-17 get corrected values for columns A,B,D,E,F,G,H,Q,R
$ws.Range("A5").Value = 111943983
$ws.Range("B5").Value = 90826
$ws.Range("E5").Value = 4366
$ws.Range("F5").Value = "Skarp dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum peckii"
$ws.Range("H5").Value = "Banker"
$ws.Range("Q5").Value = 682871
$ws.Range("R5").Value = 6694481
$ws.Range("A6").Value = 111943988
$ws.Range("B6").Value = 107576
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 220320
$ws.Range("F6").Value = "Ängsskära"
$ws.Range("G6").Value = "Serratula tinctoria"
$ws.Range("H6").Value = "L."
$ws.Range("Q6").Value = 682930
$ws.Range("R6").Value = 6694720
$ws.Range("A7").Value = 111943981
$ws.Range("B7").Value = 96640
$ws.Range("Q7").Value = 682877
$ws.Range("R7").Value = 6694410
$ws.Range("A8").Value = 111943995
$ws.Range("B8").Value = 89047
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 3286
$ws.Range("F8").Value = "Flattoppad klubbsvamp"
$ws.Range("G8").Value = "Clavariadelphus truncatus"
$ws.Range("H8").Value = "(Quél.) Donk"
$ws.Range("Q8").Value = 682779
$ws.Range("R8").Value = 6694551
$ws.Range("B9").Value = 102192
$ws.Range("A10").Value = 111943996
$ws.Range("B10").Value = 90480
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 4769
$ws.Range("F10").Value = "Svavelriska"
$ws.Range("G10").Value = "Lactarius scrobiculatus"
$ws.Range("H10").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q10").Value = 682785
$ws.Range("R10").Value = 6694547
$ws.Range("A11").Value = 111943999
$ws.Range("B11").Value = 99874
$ws.Range("Q11").Value = 682757
$ws.Range("R11").Value = 6694406
$ws.Range("A12").Value = 111943998
$ws.Range("B12").Value = 98980
$ws.Range("E12").Value = 222498
$ws.Range("F12").Value = "Blåsippa"
$ws.Range("G12").Value = "Hepatica nobilis"
$ws.Range("H12").Value = "Schreb."
$ws.Range("Q12").Value = 682757
$ws.Range("R12").Value = 6694406
$ws.Range("A13").Value = 111943997
$ws.Range("B13").Value = 96713
$ws.Range("E13").Value = 219798
$ws.Range("F13").Value = "Skogsknipprot"
$ws.Range("G13").Value = "Epipactis helleborine"
$ws.Range("H13").Value = "(L.) Crantz"
$ws.Range("Q13").Value = 682781
$ws.Range("R13").Value = 6694488
$ws.Range("A14").Value = 111943992
$ws.Range("B14").Value = 89331
$ws.Range("E14").Value = 3215
$ws.Range("F14").Value = "Rödgul trumpetsvamp"
$ws.Range("G14").Value = "Craterellus lutescens"
$ws.Range("H14").Value = "(Fr.) Fr."
$ws.Range("Q14").Value = 682867
$ws.Range("R14").Value = 6694644
$ws.Range("A15").Value = 111943980
$ws.Range("B15").Value = 89331
$ws.Range("E15").Value = 3215
$ws.Range("F15").Value = "Rödgul trumpetsvamp"
$ws.Range("G15").Value = "Craterellus lutescens"
$ws.Range("H15").Value = "(Fr.) Fr."
$ws.Range("Q15").Value = 682877
$ws.Range("R15").Value = 6694410
$ws.Range("A16").Value = 111943984
$ws.Range("B16").Value = 99874
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 221235
$ws.Range("F16").Value = "Vårärt"
$ws.Range("G16").Value = "Lathyrus vernus"
$ws.Range("H16").Value = "(L.) Bernh."
$ws.Range("Q16").Value = 682929
$ws.Range("R16").Value = 6694685
$ws.Range("A17").Value = 111943979
$ws.Range("B17").Value = 96640
$ws.Range("E17").Value = 504
$ws.Range("F17").Value = "Guckusko"
$ws.Range("G17").Value = "Cypripedium calceolus"
$ws.Range("H17").Value = "L."
$ws.Range("Q17").Value = 682879
$ws.Range("R17").Value = 6694407

# K12 -> cleared (moved to K13)
$ws.Range("K12").ClearContents()
$ws.Range("K13").Value = "i frukt"
